$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "28.887.61"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.884.72"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4589"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07844"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9857"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.79"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "1.882.36"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.682"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06934"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009956"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "28.896.19"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.273"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "2.117.19"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.086"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.31"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.994"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.928"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.43"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09332"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9044"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.277"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.264"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.200"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05762"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02071"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.002"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.639"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5666"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1766"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.689"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.248"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.93"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5358"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.849"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.533"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.069"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.13%  "
